# Generate Report for Handoff
#
# The file "7a044612-f0c2-4a94-a175-90d30bbf36a0.md" (row 3 in every
# sheet) moves from "In Translation" to "Ready for handoff", and its
# handoff timestamps are refreshed to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-18-18 20:18:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-18 20:18:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-18 20:18:17"
